$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 408, pushing existing rows 408..509 down to 409..510
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row 408 with the new record
$ws.Range("A408").Value = 3
$ws.Range("B408").Value = "Femacal de La Calera"
$ws.Range("C408").Value = "Coquimbo"
$ws.Range("D408").Value = 44754
$ws.Range("E408").Value = 5
$ws.Range("F408").Value = 100112037
$ws.Range("G408").Value = "Cebollín"
$ws.Range("H408").Value = "Sin especificar"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 250
$ws.Range("K408").Value = 6500
$ws.Range("L408").Value = 7000
$ws.Range("M408").Value = 6740
$ws.Range("N408").Value = "`$/paquete 36 unidades"
$ws.Range("O408").Value = "Provincia de Quillota"
$ws.Range("P408").Value = 187
$ws.Range("Q408").Value = 36
$ws.Range("R408").Value = "Hortaliza"
